$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; this shifts existing rows 7-77 down to 8-78,
# preserving all of their data (matches the diff's net effect of a weekly
# data refresh that prepends one new record row).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record's values.
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(7, 3).Value = "Metropolitana"
$ws.Cells.Item(7, 4).Value = 45050
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 6).Value = 100112035
$ws.Cells.Item(7, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 70
$ws.Cells.Item(7, 11).Value = 20000
$ws.Cells.Item(7, 12).Value = 22000
$ws.Cells.Item(7, 13).Value = 21000
$ws.Cells.Item(7, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(7, 16).Value = 1400
$ws.Cells.Item(7, 17).Value = 15
$ws.Cells.Item(7, 18).Value = "Hortaliza"
